# Weekly update: insert a new week's price record at the top of the data
# (row 6), pushing the existing rows 6-10 down to rows 7-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6 - shifts old rows 6..10 down to 7..11
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the new week's data
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44467
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 8500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 340
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
